{"js": "// The build removed the trailing \"\u00a9 2020 ... Creative Commons Attribution\"\n// footer paragraph together with the two blank/page-break paragraphs that\n// separated it from the \"LOB1012: Estat\u00edstica (Requisito fraco)\" line,\n// while leaving the two blank/page-break paragraphs that originally\n// followed the footer paragraph untouched.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the requirement line that the footer block directly follows.\nconst anchorText = \"LOB1012: Estat\u00edstica (Requisito fraco)\";\nlet anchor = null;\nfor (const p of paragraphs.items) {\n  if (p.text === anchorText) {\n    anchor = p;\n    break;\n  }\n}\n\nif (anchor) {\n  // Walk forward from the anchor paragraph, collecting paragraphs until\n  // (and including) the copyright/footer paragraph. Only delete them if\n  // that footer paragraph is actually found nearby - otherwise leave the\n  // document untouched.\n  const candidates = [];\n  let current = anchor.getNextOrNullObject();\n  current.load(\"text,isNullObject\");\n  await context.sync();\n\n  let footerFound = false;\n  let guard = 0;\n  while (!current.isNullObject && guard < 6) {\n    candidates.push(current);\n    if (current.text.indexOf(\"Powered by Jekyll and Github pages\") !== -1) {\n      footerFound = true;\n      break;\n    }\n    current = current.getNextOrNullObject();\n    current.load(\"text,isNullObject\");\n    await context.sync();\n    guard++;\n  }\n\n  if (footerFound) {\n    for (const p of candidates) {\n      p.delete();\n    }\n    await context.sync();\n  }\n}\n", "ps1": "# The build removed the trailing \"(c) 2020 ... Creative Commons Attribution\"\n# footer paragraph together with the two blank/page-break paragraphs that\n# separated it from the \"LOB1012: Estatistica (Requisito fraco)\" line, while\n# leaving the two blank/page-break paragraphs that originally followed the\n# footer paragraph untouched.\n\n$d = $word.ActiveDocument\n\n$targetIndex = -1\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"LOB1012*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -ge 1) {\n    # Look ahead (without mutating) for the footer/copyright paragraph so we\n    # only delete when it is actually present nearby - otherwise leave the\n    # document untouched.\n    $footerOffset = -1\n    $maxLookahead = 6\n    for ($k = 1; $k -le $maxLookahead; $k++) {\n        $idx = $targetIndex + $k\n        if ($idx -gt $d.Paragraphs.Count) {\n            break\n        }\n        $p = $d.Paragraphs.Item($idx)\n        if ($p.Range.Text -like \"*Powered by Jekyll and Github pages*\") {\n            $footerOffset = $k\n            break\n        }\n    }\n\n    if ($footerOffset -ge 1) {\n        # After each deletion the collection re-indexes, so the paragraph\n        # right after the anchor is always at $targetIndex + 1.\n        for ($n = 0; $n -lt $footerOffset; $n++) {\n            $p = $d.Paragraphs.Item($targetIndex + 1)\n            $p.Range.Delete()\n        }\n    }\n}\n"}
